# The workbook gained a new weekly price record. It is inserted as a new
# row 172 on "Sheet1", which pushes the previous rows 172-247 down to
# 173-248 (the last existing row, old 247, becomes new row 248).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 172, shifting rows 172:247 down
# to 173:248.
$ws.Rows.Item(172).Insert()

# Populate the newly inserted row 172 with the new observation.
$ws.Cells.Item(172, 1).Value = 9
$ws.Cells.Item(172, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(172, 3).Value = "Metropolitana"
$ws.Cells.Item(172, 4).Value = 44460
$ws.Cells.Item(172, 5).Value = 13
$ws.Cells.Item(172, 6).Value = 100112031
$ws.Cells.Item(172, 7).Value = "Poroto verde"
$ws.Cells.Item(172, 8).Value = "Magnum"
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 25
$ws.Cells.Item(172, 11).Value = 40000
$ws.Cells.Item(172, 12).Value = 43000
$ws.Cells.Item(172, 13).Value = 41560
$ws.Cells.Item(172, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(172, 15).Value = "Perú"
$ws.Cells.Item(172, 16).Value = 1662
$ws.Cells.Item(172, 17).Value = 25
$ws.Cells.Item(172, 18).Value = "Hortaliza"
